$wb = $excel.ActiveWorkbook

# New shared strings must land in the order: Runway, Flight, Terminal
# (matches the target sharedStrings.xml ordering), so touch the sheets
# in that order.

# --- Sheet 2: "Taxi distances" --- (adds "Runway")
$ws2 = $wb.Worksheets.Item("Taxi distances")
$ws2.Range("A1").Value = "Runway"

# --- Sheet 1: "Flight schedule" --- (adds "Flight")
$ws1 = $wb.Worksheets.Item("Flight schedule")
$ws1.Range("A1").Value = "Flight"
$ws1.Range("B15").Select()

# --- Sheet 3: "Terminal capacity" --- (adds "Terminal")
$ws3 = $wb.Worksheets.Item("Terminal capacity")
$ws3.Range("A1").Value = "Terminal"
$ws3.Range("B4").Select()
